# Auto-generated Excel COM-interop script
# Applies the exact cell-value changes described by the target diff
# (Hortaliza, Femacal de La Calera - Espinaca: weekly fruit/vegetable price update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 384
$ws.Cells.Item(384, 4).Value = 44320
$ws.Cells.Item(384, 10).Value = 80
$ws.Cells.Item(384, 11).Value = 3500
$ws.Cells.Item(384, 12).Value = 3500
$ws.Cells.Item(384, 13).Value = 3500
$ws.Cells.Item(384, 16).Value = 1167

# Row 385
$ws.Cells.Item(385, 4).Value = 44978
$ws.Cells.Item(385, 10).Value = 130
$ws.Cells.Item(385, 11).Value = 5500
$ws.Cells.Item(385, 12).Value = 6000
$ws.Cells.Item(385, 13).Value = 5731
$ws.Cells.Item(385, 16).Value = 1910

# Row 386
$ws.Cells.Item(386, 4).Value = 44819
$ws.Cells.Item(386, 10).Value = 140
$ws.Cells.Item(386, 11).Value = 4000
$ws.Cells.Item(386, 12).Value = 4500
$ws.Cells.Item(386, 13).Value = 4286
$ws.Cells.Item(386, 16).Value = 1429

# Row 387
$ws.Cells.Item(387, 4).Value = 44509
$ws.Cells.Item(387, 10).Value = 290
$ws.Cells.Item(387, 11).Value = 2000
$ws.Cells.Item(387, 12).Value = 2300
$ws.Cells.Item(387, 13).Value = 2166
$ws.Cells.Item(387, 16).Value = 722

# Row 388
$ws.Cells.Item(388, 4).Value = 44932
$ws.Cells.Item(388, 10).Value = 170
$ws.Cells.Item(388, 11).Value = 4000
$ws.Cells.Item(388, 12).Value = 4500
$ws.Cells.Item(388, 13).Value = 4235
$ws.Cells.Item(388, 16).Value = 1412

# Row 389
$ws.Cells.Item(389, 4).Value = 44622
$ws.Cells.Item(389, 10).Value = 130
$ws.Cells.Item(389, 11).Value = 5000
$ws.Cells.Item(389, 12).Value = 5500
$ws.Cells.Item(389, 13).Value = 5269
$ws.Cells.Item(389, 16).Value = 1756

# Row 390
$ws.Cells.Item(390, 4).Value = 44582
$ws.Cells.Item(390, 10).Value = 150
$ws.Cells.Item(390, 11).Value = 3500
$ws.Cells.Item(390, 12).Value = 4000
$ws.Cells.Item(390, 13).Value = 3767
$ws.Cells.Item(390, 16).Value = 1256

# Row 391
$ws.Cells.Item(391, 4).Value = 44777
$ws.Cells.Item(391, 10).Value = 230
$ws.Cells.Item(391, 11).Value = 4000
$ws.Cells.Item(391, 12).Value = 4500
$ws.Cells.Item(391, 13).Value = 4261
$ws.Cells.Item(391, 16).Value = 1420

# Row 392
$ws.Cells.Item(392, 4).Value = 45180
$ws.Cells.Item(392, 10).Value = 160
$ws.Cells.Item(392, 13).Value = 4188
$ws.Cells.Item(392, 16).Value = 1396

# Row 393
$ws.Cells.Item(393, 4).Value = 45062
$ws.Cells.Item(393, 10).Value = 90
$ws.Cells.Item(393, 12).Value = 4000
$ws.Cells.Item(393, 13).Value = 4000
$ws.Cells.Item(393, 16).Value = 1333

# Row 394
$ws.Cells.Item(394, 4).Value = 45119
$ws.Cells.Item(394, 10).Value = 65
$ws.Cells.Item(394, 11).Value = 5000
$ws.Cells.Item(394, 12).Value = 5000
$ws.Cells.Item(394, 13).Value = 5000
$ws.Cells.Item(394, 16).Value = 1667

# Row 395
$ws.Cells.Item(395, 4).Value = 44462
$ws.Cells.Item(395, 10).Value = 310
$ws.Cells.Item(395, 11).Value = 3000
$ws.Cells.Item(395, 12).Value = 3500
$ws.Cells.Item(395, 13).Value = 3242
$ws.Cells.Item(395, 16).Value = 1081

# Row 396
$ws.Cells.Item(396, 4).Value = 44285
$ws.Cells.Item(396, 10).Value = 80
$ws.Cells.Item(396, 11).Value = 4000
$ws.Cells.Item(396, 12).Value = 4000
$ws.Cells.Item(396, 13).Value = 4000
$ws.Cells.Item(396, 16).Value = 1333

# Row 397
$ws.Cells.Item(397, 4).Value = 44810
$ws.Cells.Item(397, 10).Value = 200
$ws.Cells.Item(397, 12).Value = 4500
$ws.Cells.Item(397, 13).Value = 4225
$ws.Cells.Item(397, 16).Value = 1408

# Row 398
$ws.Cells.Item(398, 4).Value = 44855
$ws.Cells.Item(398, 10).Value = 230
$ws.Cells.Item(398, 13).Value = 3739
$ws.Cells.Item(398, 16).Value = 1246

# Row 399
$ws.Cells.Item(399, 4).Value = 44599
$ws.Cells.Item(399, 10).Value = 100
$ws.Cells.Item(399, 13).Value = 3700
$ws.Cells.Item(399, 16).Value = 1233

# Row 400
$ws.Cells.Item(400, 4).Value = 45155
$ws.Cells.Item(400, 10).Value = 120
$ws.Cells.Item(400, 11).Value = 4000
$ws.Cells.Item(400, 13).Value = 4000
$ws.Cells.Item(400, 16).Value = 1333

# Row 401
$ws.Cells.Item(401, 4).Value = 44931

# Row 402
$ws.Cells.Item(402, 4).Value = 44428

# Row 403
$ws.Cells.Item(403, 4).Value = 44677
$ws.Cells.Item(403, 10).Value = 270
$ws.Cells.Item(403, 11).Value = 3500
$ws.Cells.Item(403, 13).Value = 3667
$ws.Cells.Item(403, 16).Value = 1222

# Row 404
$ws.Cells.Item(404, 4).Value = 44937
$ws.Cells.Item(404, 10).Value = 155
$ws.Cells.Item(404, 11).Value = 4500
$ws.Cells.Item(404, 12).Value = 4800
$ws.Cells.Item(404, 13).Value = 4645
$ws.Cells.Item(404, 16).Value = 1548

# Row 405
$ws.Cells.Item(405, 4).Value = 44211
$ws.Cells.Item(405, 10).Value = 60
$ws.Cells.Item(405, 11).Value = 2500
$ws.Cells.Item(405, 12).Value = 2500
$ws.Cells.Item(405, 13).Value = 2500
$ws.Cells.Item(405, 16).Value = 833

# Row 406
$ws.Cells.Item(406, 4).Value = 44291
$ws.Cells.Item(406, 10).Value = 148
$ws.Cells.Item(406, 11).Value = 3800
$ws.Cells.Item(406, 12).Value = 4000
$ws.Cells.Item(406, 13).Value = 3922
$ws.Cells.Item(406, 16).Value = 1307

# Row 407
$ws.Cells.Item(407, 4).Value = 44461
$ws.Cells.Item(407, 10).Value = 160
$ws.Cells.Item(407, 11).Value = 2500
$ws.Cells.Item(407, 12).Value = 2500
$ws.Cells.Item(407, 13).Value = 2500
$ws.Cells.Item(407, 16).Value = 833

# Row 408
$ws.Cells.Item(408, 4).Value = 45100
$ws.Cells.Item(408, 10).Value = 220
$ws.Cells.Item(408, 11).Value = 5000
$ws.Cells.Item(408, 12).Value = 5500
$ws.Cells.Item(408, 13).Value = 5250
$ws.Cells.Item(408, 16).Value = 1750

# Row 409
$ws.Cells.Item(409, 4).Value = 44726
$ws.Cells.Item(409, 10).Value = 230
$ws.Cells.Item(409, 11).Value = 3500
$ws.Cells.Item(409, 12).Value = 4000
$ws.Cells.Item(409, 13).Value = 3761
$ws.Cells.Item(409, 16).Value = 1254

# Row 410
$ws.Cells.Item(410, 4).Value = 44160
$ws.Cells.Item(410, 10).Value = 120
$ws.Cells.Item(410, 12).Value = 3500
$ws.Cells.Item(410, 13).Value = 3500
$ws.Cells.Item(410, 16).Value = 1167

# Row 411
$ws.Cells.Item(411, 4).Value = 44445
$ws.Cells.Item(411, 10).Value = 260
$ws.Cells.Item(411, 11).Value = 2500
$ws.Cells.Item(411, 12).Value = 3000
$ws.Cells.Item(411, 13).Value = 2731
$ws.Cells.Item(411, 16).Value = 910

# Row 412
$ws.Cells.Item(412, 4).Value = 44587
$ws.Cells.Item(412, 10).Value = 60
$ws.Cells.Item(412, 11).Value = 4000
$ws.Cells.Item(412, 12).Value = 4000
$ws.Cells.Item(412, 13).Value = 4000
$ws.Cells.Item(412, 16).Value = 1333

# Row 413
$ws.Cells.Item(413, 4).Value = 45075
$ws.Cells.Item(413, 10).Value = 170
$ws.Cells.Item(413, 11).Value = 5500
$ws.Cells.Item(413, 12).Value = 6000
$ws.Cells.Item(413, 13).Value = 5765
$ws.Cells.Item(413, 16).Value = 1922

# Row 414
$ws.Cells.Item(414, 4).Value = 44476
$ws.Cells.Item(414, 10).Value = 220
$ws.Cells.Item(414, 11).Value = 2500
$ws.Cells.Item(414, 12).Value = 2500
$ws.Cells.Item(414, 13).Value = 2500
$ws.Cells.Item(414, 16).Value = 833

# Row 415
$ws.Cells.Item(415, 4).Value = 44720
$ws.Cells.Item(415, 10).Value = 210
$ws.Cells.Item(415, 11).Value = 3500
$ws.Cells.Item(415, 12).Value = 3800
$ws.Cells.Item(415, 13).Value = 3657
$ws.Cells.Item(415, 16).Value = 1219

# Row 416
$ws.Cells.Item(416, 4).Value = 45202
$ws.Cells.Item(416, 10).Value = 90
$ws.Cells.Item(416, 11).Value = 4000
$ws.Cells.Item(416, 12).Value = 4000
$ws.Cells.Item(416, 13).Value = 4000
$ws.Cells.Item(416, 16).Value = 1333
